$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "release/8.0.10"
$ws.Range("B13").Value = "X"
$ws.Range("C13").Value = "X"
$ws.Range("D13").Value = "X"
$ws.Range("E13").Value = "X"
